$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the leftover "Hyperlink" cell style from A2 (hyperlink relationship itself is kept)
$ws.Range("A2").ClearFormats()

# Row 3: new credential pair, each one a real hyperlink (mailto:) like the existing A2 entry
$ws.Range("A3").Value = "zahira@credosystemz.sandbox"
$ws.Range("B3").Value = "Waseem@20"

$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:zahira@credosystemz.sandbox")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:Waseem@20")

# Re-apply the workbook's built-in Hyperlink style so the new cells render/save
# with the same cellXf index (1) that A2 used to use.
$ws.Range("A3").Style = "Hyperlink"
$ws.Range("B3").Style = "Hyperlink"

# Match the recorded selection left behind on the sheet view
[void]$ws.Range("B2").Select()
